{"js": "// Replace each two-digit-division answer cell's text with the updated value.\n// Each old answer string is unique within the document, so a direct\n// search-and-replace on the document body safely targets the correct cell.\nconst replacements = [\n  [\"29\u00f74=7, 1\", \"44\u00f73=14, 2\"],\n  [\"10\u00f77=1, 3\", \"31\u00f76=5, 1\"],\n  [\"55\u00f77=7, 6\", \"42\u00f75=8, 2\"],\n  [\"40\u00f78=5, 0\", \"13\u00f79=1, 4\"],\n  [\"42\u00f77=6, 0\", \"78\u00f79=8, 6\"],\n  [\"92\u00f74=23, 0\", \"39\u00f74=9, 3\"],\n  [\"49\u00f72=24, 1\", \"85\u00f78=10, 5\"],\n  [\"51\u00f76=8, 3\", \"20\u00f79=2, 2\"],\n  [\"24\u00f79=2, 6\", \"38\u00f78=4, 6\"],\n  [\"22\u00f72=11, 0\", \"61\u00f77=8, 5\"],\n  [\"64\u00f73=21, 1\", \"89\u00f77=12, 5\"],\n  [\"82\u00f78=10, 2\", \"96\u00f74=24, 0\"],\n  [\"51\u00f77=7, 2\", \"74\u00f77=10, 4\"],\n  [\"54\u00f78=6, 6\", \"78\u00f73=26, 0\"],\n  [\"90\u00f75=18, 0\", \"32\u00f77=4, 4\"],\n  [\"66\u00f74=16, 2\", \"52\u00f75=10, 2\"],\n  [\"25\u00f75=5, 0\", \"54\u00f75=10, 4\"],\n  [\"36\u00f78=4, 4\", \"44\u00f73=14, 2\"],\n  [\"49\u00f73=16, 1\", \"82\u00f76=13, 4\"],\n  [\"78\u00f75=15, 3\", \"16\u00f72=8, 0\"],\n  [\"97\u00f77=13, 6\", \"15\u00f79=1, 6\"],\n  [\"35\u00f75=7, 0\", \"55\u00f76=9, 1\"],\n  [\"65\u00f76=10, 5\", \"78\u00f78=9, 6\"],\n  [\"19\u00f76=3, 1\", \"37\u00f79=4, 1\"],\n  [\"85\u00f74=21, 1\", \"77\u00f75=15, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "# Replace each two-digit-division answer cell's text with the updated value.\n# Each old answer string is unique within the document, so Find/Replace\n# targeting the whole document content safely updates the correct cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"29\u00f74=7, 1\", \"44\u00f73=14, 2\"),\n    @(\"10\u00f77=1, 3\", \"31\u00f76=5, 1\"),\n    @(\"55\u00f77=7, 6\", \"42\u00f75=8, 2\"),\n    @(\"40\u00f78=5, 0\", \"13\u00f79=1, 4\"),\n    @(\"42\u00f77=6, 0\", \"78\u00f79=8, 6\"),\n    @(\"92\u00f74=23, 0\", \"39\u00f74=9, 3\"),\n    @(\"49\u00f72=24, 1\", \"85\u00f78=10, 5\"),\n    @(\"51\u00f76=8, 3\", \"20\u00f79=2, 2\"),\n    @(\"24\u00f79=2, 6\", \"38\u00f78=4, 6\"),\n    @(\"22\u00f72=11, 0\", \"61\u00f77=8, 5\"),\n    @(\"64\u00f73=21, 1\", \"89\u00f77=12, 5\"),\n    @(\"82\u00f78=10, 2\", \"96\u00f74=24, 0\"),\n    @(\"51\u00f77=7, 2\", \"74\u00f77=10, 4\"),\n    @(\"54\u00f78=6, 6\", \"78\u00f73=26, 0\"),\n    @(\"90\u00f75=18, 0\", \"32\u00f77=4, 4\"),\n    @(\"66\u00f74=16, 2\", \"52\u00f75=10, 2\"),\n    @(\"25\u00f75=5, 0\", \"54\u00f75=10, 4\"),\n    @(\"36\u00f78=4, 4\", \"44\u00f73=14, 2\"),\n    @(\"49\u00f73=16, 1\", \"82\u00f76=13, 4\"),\n    @(\"78\u00f75=15, 3\", \"16\u00f72=8, 0\"),\n    @(\"97\u00f77=13, 6\", \"15\u00f79=1, 6\"),\n    @(\"35\u00f75=7, 0\", \"55\u00f76=9, 1\"),\n    @(\"65\u00f76=10, 5\", \"78\u00f78=9, 6\"),\n    @(\"19\u00f76=3, 1\", \"37\u00f79=4, 1\"),\n    @(\"85\u00f74=21, 1\", \"77\u00f75=15, 2\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
